$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7399.4
$ws.Range("I19").Value = 2000
$ws.Range("J19").Value = 8749.25
$ws.Range("K19").Value = 2000
$ws.Range("L19").Value = 8749.25
$ws.Range("M19").Value = -1825
$ws.Range("N19").Value = -9099.25

$ws.Range("H38").Value = 4808.353
$ws.Range("J38").Value = 6547.3
$ws.Range("L38").Value = 19641.9
$ws.Range("N38").Value = -20385.9

$ws.Range("H55").Value = 84.57143000000001
$ws.Range("I55").Value = 40
$ws.Range("J55").Value = 118
$ws.Range("K55").Value = 40
$ws.Range("L55").Value = 118
$ws.Range("M55").Value = 174
$ws.Range("N55").Value = -546

$ws.Range("H124").Value = 25200
$ws.Range("I124").Value = 25200
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 25200
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -20290
$ws.Range("N124").Value = $null

$ws.Range("H125").Value = 8637
$ws.Range("I125").Value = 7999.6665
$ws.Range("J125").Value = 9019.4
$ws.Range("K125").Value = 71996.9985
$ws.Range("L125").Value = 81174.59999999999
$ws.Range("M125").Value = -69536.9985
$ws.Range("N125").Value = -86094.59999999999

$ws.Range("H127").Value = 1572.8
$ws.Range("I127").Value = 799
$ws.Range("J127").Value = 2733.5
$ws.Range("K127").Value = 2397
$ws.Range("L127").Value = 8200.5
$ws.Range("M127").Value = 2563
$ws.Range("N127").Value = -18120.5

$ws.Range("H130").Value = 97333
$ws.Range("J130").Value = 95999.5
$ws.Range("L130").Value = 95999.5
$ws.Range("N130").Value = -106039.5

$ws.Range("J134").Value = 96250
$ws.Range("L134").Value = 96250
$ws.Range("N134").Value = -106390

$ws.Range("H135").Value = 437
$ws.Range("I135").Value = 407.96
$ws.Range("K135").Value = 3671.64
$ws.Range("M135").Value = -1136.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5596.4
$ws.Range("I2").Value = 4002.6667
$ws.Range("J2").Value = 7987
$ws.Range("K2").Value = 4002.6667
$ws.Range("L2").Value = 7987
$ws.Range("M2").Value = -3889.6667
$ws.Range("N2").Value = -8213

$ws.Range("H12").Value = 2833.3333
$ws.Range("J12").Value = 2833.3333
$ws.Range("L12").Value = 2833.3333
$ws.Range("N12").Value = -3179.3333

$ws.Range("H74").Value = 5704.2334
$ws.Range("I74").Value = 5540.25
$ws.Range("K74").Value = 5540.25
$ws.Range("M74").Value = -4666.25

$ws.Range("H77").Value = 5704.2334
$ws.Range("I77").Value = 5540.25
$ws.Range("K77").Value = 27701.25
$ws.Range("M77").Value = -23333.25

$ws.Range("H97").Value = 1449.24
$ws.Range("J97").Value = 1444.8334
$ws.Range("L97").Value = 1444.8334
$ws.Range("N97").Value = -2436.8334

$ws.Range("H116").Value = 5596.4
$ws.Range("I116").Value = 4002.6667
$ws.Range("J116").Value = 7987
$ws.Range("K116").Value = 4002.6667
$ws.Range("L116").Value = 7987
$ws.Range("M116").Value = -1708.6667
$ws.Range("N116").Value = -12575

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5596.4
$ws.Range("I3").Value = 4002.6667
$ws.Range("J3").Value = 7987
$ws.Range("K3").Value = 4002.6667
$ws.Range("L3").Value = 7987
$ws.Range("M3").Value = -3888.6667
$ws.Range("N3").Value = -8215

$ws.Range("H20").Value = 2179
$ws.Range("I20").Value = 1746.6
$ws.Range("J20").Value = 2899.6667
$ws.Range("K20").Value = 1746.6
$ws.Range("L20").Value = 2899.6667
$ws.Range("M20").Value = -1499.6
$ws.Range("N20").Value = -3393.6667

$ws.Range("H86").Value = 14916.5625
$ws.Range("I86").Value = 2369.3333
$ws.Range("J86").Value = 31048.715
$ws.Range("K86").Value = 2369.3333
$ws.Range("L86").Value = 31048.715
$ws.Range("M86").Value = -1246.3333
$ws.Range("N86").Value = -33294.715

$ws.Range("H89").Value = 14916.5625
$ws.Range("I89").Value = 2369.3333
$ws.Range("J89").Value = 31048.715
$ws.Range("K89").Value = 11846.6665
$ws.Range("L89").Value = 155243.575
$ws.Range("M89").Value = -6230.666499999999
$ws.Range("N89").Value = -166475.575

$ws.Range("H94").Value = 1281.0333
$ws.Range("I94").Value = 823.8946999999999
$ws.Range("J94").Value = 2070.6365
$ws.Range("K94").Value = 823.8946999999999
$ws.Range("L94").Value = 2070.6365
$ws.Range("M94").Value = -372.8946999999999
$ws.Range("N94").Value = -2972.6365

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 103.71429
$ws.Range("I7").Value = 45.25
$ws.Range("J7").Value = 139.6923
$ws.Range("K7").Value = 45.25
$ws.Range("L7").Value = 139.6923
$ws.Range("M7").Value = 67.75
$ws.Range("N7").Value = -365.6923

$ws.Range("H12").Value = 1287.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1287.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1287.5
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = -1627.5

$ws.Range("H94").Value = 2471.2778
$ws.Range("I94").Value = 1949.75
$ws.Range("J94").Value = 2620.2856
$ws.Range("K94").Value = 1949.75
$ws.Range("L94").Value = 2620.2856
$ws.Range("M94").Value = -1498.75
$ws.Range("N94").Value = -3522.2856

$ws.Range("H99").Value = 10216.1
$ws.Range("I99").Value = 6500
$ws.Range("J99").Value = 11145.125
$ws.Range("K99").Value = 6500
$ws.Range("L99").Value = 11145.125
$ws.Range("M99").Value = -5002
$ws.Range("N99").Value = -14141.125

$ws.Range("H123").Value = 39997.5
$ws.Range("I123").Value = 30000
$ws.Range("J123").Value = 49995
$ws.Range("K123").Value = 30000
$ws.Range("L123").Value = 49995
$ws.Range("M123").Value = -25100
$ws.Range("N123").Value = -59795

$ws.Range("H126").Value = 10216.1
$ws.Range("I126").Value = 6500
$ws.Range("J126").Value = 11145.125
$ws.Range("K126").Value = 19500
$ws.Range("L126").Value = 33435.375
$ws.Range("M126").Value = -17030
$ws.Range("N126").Value = -38375.375

$ws.Range("H134").Value = 3598.4443
$ws.Range("I134").Value = 3295.5454
$ws.Range("J134").Value = 4074.4285
$ws.Range("K134").Value = 9886.636200000001
$ws.Range("L134").Value = 12223.2855
$ws.Range("M134").Value = -7351.636200000001
$ws.Range("N134").Value = -17293.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 482.5
$ws.Range("I3").Value = 342.2857
$ws.Range("J3").Value = 809.6667
$ws.Range("K3").Value = 342.2857
$ws.Range("L3").Value = 809.6667
$ws.Range("M3").Value = -226.2857
$ws.Range("N3").Value = -1041.6667

$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 800
$ws.Range("K17").Value = 200
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = -32
$ws.Range("N17").Value = -1136

$ws.Range("H29").Value = 24000
$ws.Range("I29").Value = 24000
$ws.Range("K29").Value = 24000
$ws.Range("M29").Value = -23710

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4964
$ws.Range("I7").Value = 4019.7693
$ws.Range("J7").Value = 6498.375
$ws.Range("K7").Value = 4019.7693
$ws.Range("L7").Value = 6498.375
$ws.Range("M7").Value = -3907.7693
$ws.Range("N7").Value = -6722.375

$ws.Range("H16").Value = 1423.1936
$ws.Range("I16").Value = 1254.9584
$ws.Range("K16").Value = 1254.9584
$ws.Range("M16").Value = -1084.9584

$ws.Range("H40").Value = 3902.2942
$ws.Range("I40").Value = 3233.7
$ws.Range("K40").Value = 3233.7
$ws.Range("M40").Value = -3097.7

$ws.Range("H82").Value = 4535.769
$ws.Range("J82").Value = 3999.4
$ws.Range("L82").Value = 3999.4
$ws.Range("N82").Value = -4721.4

$ws.Range("H85").Value = 4535.769
$ws.Range("J85").Value = 3999.4
$ws.Range("L85").Value = 3999.4
$ws.Range("N85").Value = -6495.4

$ws.Range("H100").Value = 3680.6316
$ws.Range("I100").Value = 3239.375
$ws.Range("K100").Value = 3239.375
$ws.Range("M100").Value = -2698.375

$ws.Range("H126").Value = 4964
$ws.Range("I126").Value = 4019.7693
$ws.Range("J126").Value = 6498.375
$ws.Range("K126").Value = 12059.3079
$ws.Range("L126").Value = 19495.125
$ws.Range("M126").Value = -9589.3079
$ws.Range("N126").Value = -24435.125

$ws.Range("H132").Value = 6666.2705
$ws.Range("I132").Value = 1475.1613
$ws.Range("K132").Value = 4425.4839
$ws.Range("M132").Value = -1895.4839

$ws.Range("H136").Value = 4462.6924
$ws.Range("I136").Value = 4552.0835
$ws.Range("K136").Value = 13656.2505
$ws.Range("M136").Value = -11106.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8112.375
$ws.Range("I62").Value = 6499
$ws.Range("K62").Value = 6499
$ws.Range("M62").Value = -5875

$ws.Range("H65").Value = 8112.375
$ws.Range("I65").Value = 6499
$ws.Range("K65").Value = 32495
$ws.Range("M65").Value = -29375

$ws.Range("H96").Value = 4862.6113
$ws.Range("I96").Value = 3994.2222
$ws.Range("J96").Value = 5731
$ws.Range("K96").Value = 3994.2222
$ws.Range("L96").Value = 5731
$ws.Range("M96").Value = -2621.2222
$ws.Range("N96").Value = -8477

$ws.Range("H122").Value = 4918.8
$ws.Range("I122").Value = 4918.8
$ws.Range("K122").Value = 14756.4
$ws.Range("M122").Value = -12306.4
